$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 386
$ws.Range("F4").Value = 436
$ws.Range("F5").Value = 1178
$ws.Range("F8").Value = 1196
$ws.Range("F9").Value = 1670
$ws.Range("F10").Value = 6182
$ws.Range("F12").Value = 1799
$ws.Range("F13").Value = 471
$ws.Range("F15").Value = 9
$ws.Range("F18").Value = 11
$ws.Range("F19").Value = 6492
$ws.Range("F22").Value = 163
$ws.Range("F24").Value = 1697
$ws.Range("F25").Value = 840
$ws.Range("F26").Value = 11
$ws.Range("F28").Value = 158
$ws.Range("F29").Value = 1528
$ws.Range("F30").Value = 749
$ws.Range("F31").Value = 298
$ws.Range("F34").Value = 43

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 20
$ws.Range("F4").Value = 329
$ws.Range("F5").Value = 193
$ws.Range("F8").Value = 423

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9525
$ws.Range("F3").Value = 2250
$ws.Range("F5").Value = 232

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9525
$ws.Range("F3").Value = 2250
$ws.Range("F5").Value = 386
$ws.Range("F6").Value = 436
$ws.Range("F7").Value = 1178
$ws.Range("F10").Value = 20
$ws.Range("F11").Value = 329
$ws.Range("F12").Value = 1196
$ws.Range("F13").Value = 232
$ws.Range("F14").Value = 1670
$ws.Range("F15").Value = 6182
$ws.Range("F17").Value = 1799
$ws.Range("F19").Value = 471
$ws.Range("F21").Value = 9
$ws.Range("F23").Value = 6492
$ws.Range("F26").Value = 163
$ws.Range("F28").Value = 1697
$ws.Range("F29").Value = 840
$ws.Range("F31").Value = 158
$ws.Range("F32").Value = 1528
$ws.Range("F33").Value = 749
$ws.Range("F35").Value = 298
